$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "会议主题" / ${record.subject} column (column I) is removed from the
# DepartmentMeetings template. Deleting the entire column shifts the
# following columns (J/K/L -> I/J/K) left and drops the now-unused
# "会议主题" / "${record.subject}" entries from the shared string table.
$ws.Columns("I").Delete()

# Match the author's resulting selection (column I, now "费用预算" /
# ${record.planCost}) left active after the delete.
$ws.Columns("I").Select() | Out-Null
